# Weekly update: two new price records for "Zanahoria" (Carrot) are
# published, so they get inserted at the top of the data block (row 262),
# pushing all the previously-existing rows (old 262-282) down by two rows
# (new 264-284).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at row 262; everything below (old rows
# 262-282) shifts down to become rows 264-284.
$ws.Rows.Item(262).Insert()
$ws.Rows.Item(262).Insert()

# --- New row 262 ---
$ws.Cells.Item(262, 1).Value  = 1
$ws.Cells.Item(262, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(262, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(262, 4).Value  = 44714
$ws.Cells.Item(262, 5).Value  = 15
$ws.Cells.Item(262, 6).Value  = 100114013
$ws.Cells.Item(262, 7).Value  = "Zanahoria"
$ws.Cells.Item(262, 8).Value  = "Sin especificar"
$ws.Cells.Item(262, 9).Value  = "Primera"
$ws.Cells.Item(262, 10).Value = 35
$ws.Cells.Item(262, 11).Value = 17000
$ws.Cells.Item(262, 12).Value = 17000
$ws.Cells.Item(262, 13).Value = 17000
$ws.Cells.Item(262, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(262, 15).Value = "Provincia de Calama"
$ws.Cells.Item(262, 16).Value = 680
$ws.Cells.Item(262, 17).Value = 25
$ws.Cells.Item(262, 18).Value = "Hortaliza"

# --- New row 263 ---
$ws.Cells.Item(263, 1).Value  = 1
$ws.Cells.Item(263, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(263, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(263, 4).Value  = 44714
$ws.Cells.Item(263, 5).Value  = 15
$ws.Cells.Item(263, 6).Value  = 100114013
$ws.Cells.Item(263, 7).Value  = "Zanahoria"
$ws.Cells.Item(263, 8).Value  = "Sin especificar"
$ws.Cells.Item(263, 9).Value  = "Segunda"
$ws.Cells.Item(263, 10).Value = 35
$ws.Cells.Item(263, 11).Value = 16000
$ws.Cells.Item(263, 12).Value = 16000
$ws.Cells.Item(263, 13).Value = 16000
$ws.Cells.Item(263, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(263, 15).Value = "Provincia de Calama"
$ws.Cells.Item(263, 16).Value = 640
$ws.Cells.Item(263, 17).Value = 25
$ws.Cells.Item(263, 18).Value = "Hortaliza"
